$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table (values scraped on
# Wed Aug  7 21:47:21 UTC 2024), including the Stacks/Filecoin (rows 42-43)
# and EnergySwap/Maker (rows 50-51) row-order swaps.

$ws.Range("D2").Value = "55.413.43"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "2.354.30"
$ws.Range("E3").Value = "  -4.94%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'476.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").Value = "'146.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +23.79%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "2.356.39"
$ws.Range("E9").Value = "  -5.19%  "
$ws.Range("D10").Value = "'0.0968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "'5.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.12%  "
$ws.Range("D12").Value = "'0.326"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "2.753.16"
$ws.Range("E14").Value = "  -5.56%  "
$ws.Range("D15").Value = "55.247.13"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "'20.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("E17").Value = "  -4.84%  "
$ws.Range("D18").Value = "2.346.54"
$ws.Range("E18").Value = "  -5.71%  "
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "'315.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'9.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.66%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'5.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").Value = "'57.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'0.395"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").Value = "'0.152"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.38%  "
$ws.Range("D28").Value = "2.445.78"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").Value = "'7.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.04%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "0.0₃0747"
$ws.Range("E31").Value = "  -5.26%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "'145.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("D38").Value = "'0.810"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.75%  "
$ws.Range("E39").Value = "  +9.90%  "
$ws.Range("D40").Value = "'33.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "'0.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").Value = "'0.577"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("D45").Value = "'0.0518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.59%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'252.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").Value = "'0.0220"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").Value = "'4.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.21%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.790.57"
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'16.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.95%  "
